# Apply the "Added new test case and selenium grid architecture" edit to the
# Login sheet of the TestData workbook.
#
# Summary of the content change:
#   - Row 2 (the "Login" test case) now uses a new mailbox/password pair:
#       B2: testenterprise01@mailinator.com -> qa_testers@qabrains.com
#       C2: Test@1234                       -> Password123
#   - Row 3 becomes a new "Forgot Password" test case:
#       A3: Account Details -> Forgot Password
#       B3: keeps testenterprise01@mailinator.com (still a mailto hyperlink)
#       C3: the password value/hyperlink is removed (cell left blank, but
#           keeps the Hyperlink cell style)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Remove all existing hyperlinks (and their relationships) up front so we can
# cleanly rebuild only the ones that remain in the new layout. Recreate them
# in B3, B2, C2 order (their final relative order) before touching any cell
# values, since adding a hyperlink on its own does not disturb cell text.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:testenterprise01@mailinator.com")
$ws.Range("B3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:testenterprise01@mailinator.com")
$ws.Range("B2").Style = "Hyperlink"

# C2 keeps pointing at the original "Test@1234" mailto address/display text
# even though the cell itself will now show the new password.
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Test@1234", [Type]::Missing, [Type]::Missing, "Test@1234")
$ws.Range("C2").Style = "Hyperlink"

# C3 no longer carries a password / hyperlink, but still keeps the Hyperlink
# cell style that the column previously had.
$ws.Range("C3").Style = "Hyperlink"

# Now update the cell text. B2 and C2 first (so the freed-up "Test@1234" /
# "Account Details" shared strings are dropped before the new "Forgot
# Password" string is introduced), then clear C3, then rename A3.
$ws.Range("B2").Value = "qa_testers@qabrains.com"
$ws.Range("C2").Value = "Password123"
$ws.Range("C3").Value = ""
$ws.Range("A3").Value = "Forgot Password"
$ws.Range("B3").Value = "testenterprise01@mailinator.com"

# Reflect the cell that was selected when the workbook was last saved.
$ws.Range("B3").Select()
